$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time values backing the literal-formula validations (A1, B1).
# 0.0416666666666667 (~1:00:00) -> 0.25 (6:00:00)
# 0.104166666666667  (~2:30:00) -> 0.5  (12:00:00)
$ws.Range("A1").Value = 0.25
$ws.Range("B1").Value = 0.5

# The first set of data validations (rows 2-10) use literal formula values
# mirroring A1/B1 (rather than a cell reference like the second set, rows
# 11-20, which uses $A$1/$B$1). Refresh those literals to match the new
# A1/B1 values, preserving each range's type/operator/alert style.
$singleFormulaRanges = @("A2:A10", "B2:B10", "C2:C10", "D2:D10", "E2:E10", "F2:F10")
foreach ($addr in $singleFormulaRanges) {
    $v = $ws.Range($addr).Validation
    $v.Modify($v.Type, $v.AlertStyle, $v.Operator, 0.25)
}

$betweenRanges = @("G2:G10", "H2:H10")
foreach ($addr in $betweenRanges) {
    $v = $ws.Range($addr).Validation
    $v.Modify($v.Type, $v.AlertStyle, $v.Operator, 0.25, 0.5)
}
